# Updated: - risk influence - took out supply chain
# => take this version for interim presentation
#
# Data/uncertain_variables.xlsx: split "uncert_tree_vulnerability" into a
# mean/var pair, fix the "wholesail_threshhold" typos, and add a new
# "uncert_quali_threshold" variable row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: uncert_tree_vulnerability -> uncert_tree_vulnerability_mean (upper 0.9 -> 0.3) ---
$ws.Range("A2").Value = "uncert_tree_vulnerability_mean"
$ws.Range("B2").Value = "tnorm_0_1"
$ws.Range("C2").Value = 0.1
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = 0.3
$ws.Range("F2").Value = "absolute"
$ws.Range("G2").Value = "Unknown factors which make the trees vulnerable"

# --- Row 3 (new): uncert_tree_vulnerability_var, pushes the rest down ---
$ws.Range("A3").Value = "uncert_tree_vulnerability_var"
$ws.Range("B3").Value = "tnorm_0_1"
$ws.Range("C3").Value = 0.1
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = 0.6
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = $null

# --- Row 4: uncert_tree_parameter_age_1 (unchanged values, shifted from row 3) ---
$ws.Range("A4").Value = "uncert_tree_parameter_age_1"
$ws.Range("B4").Value = "posnorm"
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = """Best ages"", Verschiebung auf y-Achse"

# --- Row 5: uncert_tree_parameter_age_2 (lower/upper 7/9 -> 5/7, shifted from row 4) ---
$ws.Range("A5").Value = "uncert_tree_parameter_age_2"
$ws.Range("B5").Value = "posnorm"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = """Kurvenstauchung"""

# --- Row 6: uncert_wholesail_threshhold -> uncert_wholesail_threshold (typo fix, shifted from row 5) ---
$ws.Range("A6").Value = "uncert_wholesail_threshold"
$ws.Range("B6").Value = "posnorm"
$ws.Range("C6").Value = 4000
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = 5000
$ws.Range("F6").Value = "kg"
$ws.Range("G6").Value = "If ""too much fruits"" are produced, it must be selled to wholesome market"

# --- Row 7: uncert_influence_quali (unchanged, shifted from row 6) ---
$ws.Range("A7").Value = "uncert_influence_quali"
$ws.Range("B7").Value = "tnorm_0_1"
$ws.Range("C7").Value = 0.1
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = 0.9
$ws.Range("F7").Value = "percent per Eur"
$ws.Range("G7").Value = "Unknown factor that shows the influence of fruit quality"

# --- Row 8: uncert_influence_supply_chain_invest (unchanged, shifted from row 7) ---
$ws.Range("A8").Value = "uncert_influence_supply_chain_invest"
$ws.Range("B8").Value = "tnorm_0_1"
$ws.Range("C8").Value = 0.1
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = 0.9
$ws.Range("F8").Value = "percent per Eur"
$ws.Range("G8").Value = "Unknown factor that shows the influence of supply chain building"

# --- Row 9: uncert_risk_decrease_mean (unchanged, shifted from row 8, gains description) ---
$ws.Range("A9").Value = "uncert_risk_decrease_mean"
$ws.Range("B9").Value = "tnorm_0_1"
$ws.Range("C9").Value = 0.1
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = 0.9
$ws.Range("F9").Value = "absolute"
$ws.Range("G9").Value = "Unknown influence of yield reliability"

# --- Row 10: uncert_risk_decrease_var (unchanged, shifted from row 9) ---
$ws.Range("A10").Value = "uncert_risk_decrease_var"
$ws.Range("B10").Value = "tnorm_0_1"
$ws.Range("C10").Value = 0.1
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = 0.9
$ws.Range("F10").Value = "absolute"
$ws.Range("G10").Value = $null

# --- Row 11: uncert_wholesail_threshhold_t -> uncert_wholesail_threshold_t (typo fix, shifted from row 10) ---
$ws.Range("A11").Value = "uncert_wholesail_threshold_t"
$ws.Range("B11").Value = "posnorm"
$ws.Range("C11").Value = 1000
$ws.Range("C11").NumberFormat = "0.00"
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = 1500
$ws.Range("F11").Value = "kg"
$ws.Range("G11").Value = $null

# --- Row 12 (new): uncert_quali_threshold ---
$ws.Range("A12").Value = "uncert_quali_threshold"
$ws.Range("B12").Value = "tnorm_0_1"
$ws.Range("C12").Value = 0.2
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = 0.7
$ws.Range("F12").Value = "absolute"
$ws.Range("G12").Value = "minimum quality for direct marketing"

# --- Row 15 (new blank formatted row at the end) ---
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

# Match the saved selection/view state from the diff (E13, no frozen topLeftCell)
$ws.Range("E13").Select()
